$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header (F1 = "Trening"), matching the formatting used by the
# existing header cells (bold font, border, centered/top alignment).
$ws.Range("F1").Value = "Trening"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

function Set-Row($r, $a, $b, $c, $d, $e, $f) {
    $row = New-Object 'object[,]' 1,6
    $row[0,0] = $a
    $row[0,1] = $b
    $row[0,2] = $c
    $row[0,3] = $d
    $row[0,4] = $e
    $row[0,5] = $f
    $ws.Range("A" + $r + ":F" + $r).Value = $row
}

Set-Row 2  45684.5914462963   500.9   11.82 1.903201307569232 "10-15" "Duża Gra"
Set-Row 3  45684.59210486111  557.8   11.38 1.891998989241466 "10-15" "Duża Gra"
Set-Row 4  45684.59272986111  611.8   12.64 1.631305115563529 "10-15" "Duża Gra"
Set-Row 5  45684.59144282407  500.6   9.52  1.685560260500227 "5-10"  "Duża Gra"
Set-Row 6  45684.5921025463   557.6   9.69  1.795107943671091 "5-10"  "Duża Gra"
Set-Row 7  45684.59337916667  667.9   9.65  1.530354601996286 "5-10"  "Duża Gra"
Set-Row 8  45684.59850416666  1110.7  13.95 2.984936680112565 "10-15" "Mała Gra"
Set-Row 9  45684.60123564815  1346.7  13.83 3.271494235311233 "10-15" "Mała Gra"
Set-Row 10 45684.60197986111  1411    14.11 3.578834329332624 "10-15" "Mała Gra"
Set-Row 11 45684.59683055556  966.1   8.710000000000001 2.717948845454625 "5-10" "Mała Gra"
Set-Row 12 45684.59755393519  1028.6  9.09  2.593538182122364 "5-10"  "Mała Gra"
Set-Row 13 45684.6026337963   1467.5  9.44  2.481043015207563 "5-10"  "Mała Gra"

# Apply the date/time number format to column A (rows 2-13), matching the
# author's workflow of trying a lowercase format on one cell first, then
# settling on an uppercase format applied to the whole range.
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2:A13").NumberFormat = "YYYY-MM-DD HH:MM:SS"
